# Added check/logic for players in year 2 or 3 on new contract to only
# have vet min 1 year deal.
#
# This updates the "Rating Range Start" (col B) and "Rating Range End"
# (col C) thresholds on the ExpectedContractLength sheet for the
# 2-year and 3-year (and a few other) contract-length rows, per position,
# to reflect the new vet-minimum-1-year-deal rule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# QB
$ws.Range("B2").Value = 85
$ws.Range("B3").Value = 76
$ws.Range("C3").Value = 84

# FB
$ws.Range("B11").Value = 82
$ws.Range("C12").Value = 81

# WR
$ws.Range("B14").Value = 95
$ws.Range("C15").Value = 94
$ws.Range("B16").Value = 78
$ws.Range("C17").Value = 77

# TE
$ws.Range("B19").Value = 85
$ws.Range("B20").Value = 74
$ws.Range("C20").Value = 84
$ws.Range("B21").Value = 68
$ws.Range("C21").Value = 73
$ws.Range("C22").Value = 67

# LT
$ws.Range("B24").Value = 81
$ws.Range("C25").Value = 80

# LG
$ws.Range("B28").Value = 94
$ws.Range("B29").Value = 81
$ws.Range("C29").Value = 93
$ws.Range("B30").Value = 74
$ws.Range("C30").Value = 80
$ws.Range("C31").Value = 73

# C
$ws.Range("B33").Value = 95
$ws.Range("B34").Value = 81
$ws.Range("C34").Value = 94
$ws.Range("B35").Value = 75
$ws.Range("C35").Value = 80
$ws.Range("B36").Value = 70
$ws.Range("C36").Value = 74
$ws.Range("C37").Value = 69

# RG
$ws.Range("B38").Value = 94
$ws.Range("B39").Value = 81
$ws.Range("C39").Value = 93
$ws.Range("B40").Value = 74
$ws.Range("C40").Value = 80
$ws.Range("C41").Value = 73

# RT
$ws.Range("B44").Value = 81
$ws.Range("C45").Value = 80
$ws.Range("B46").Value = 68
$ws.Range("C47").Value = 67

# LE
$ws.Range("B48").Value = 95
$ws.Range("C49").Value = 94
$ws.Range("B50").Value = 75
$ws.Range("C51").Value = 74

# RE
$ws.Range("B53").Value = 95
$ws.Range("C54").Value = 94
$ws.Range("B55").Value = 75
$ws.Range("C56").Value = 74

# DT
$ws.Range("B60").Value = 75
$ws.Range("C61").Value = 74

# LOLB
$ws.Range("B63").Value = 90
$ws.Range("C64").Value = 89

# MLB
$ws.Range("B67").Value = 95
$ws.Range("C68").Value = 94
$ws.Range("B70").Value = 72
$ws.Range("C71").Value = 71

# ROLB
$ws.Range("B72").Value = 90
$ws.Range("C73").Value = 89

# CB
$ws.Range("B76").Value = 94
$ws.Range("C77").Value = 93

# FS
$ws.Range("B83").Value = 71
$ws.Range("C84").Value = 70

# SS
$ws.Range("B87").Value = 71
$ws.Range("C88").Value = 70

# K
$ws.Range("B91").Value = 75
$ws.Range("C92").Value = 74

# Restore the cursor/selection to match the author's final editing
# position (bottom of the sheet).
$ws.Range("B97").Select()
